# Update "countries & provincias Spain" COVID data sheet.
# The sheet ("Pais") lists countries sorted descending by total cases
# (column B). This refresh brings in newer per-country figures; a few
# countries (Paraguay, Burkina Faso, Santa Lucia) grew enough to overtake
# their neighbours in the ranking, so those rows swap names along with
# their new values. The footer timestamp is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer: "Datos actualizados a ..." timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 02:53"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 7361206
$ws.Range("C4").Value = 37013
$ws.Range("D4").Value = 4606965
$ws.Range("E4").Value = 2544464
$ws.Range("G4").Value = 324
$ws.Range("H4").Value = 209777

# Row 12: Argentina - refreshed totals
$ws.Range("B12").Value = 723132
$ws.Range("C12").Value = 11807
$ws.Range("D12").Value = 576715
$ws.Range("E12").Value = 130304
$ws.Range("G12").Value = 364
$ws.Range("H12").Value = 16113

# Rows 69-71: Paraguay's new total count overtakes Afganistan and
# Estado de Palestina, so it moves up to row 69; the other two shift
# down a row each (their own figures are unchanged).
$ws.Range("A69").Value = "Paraguay"
$ws.Range("B69").Value = 39432
$ws.Range("C69").Value = 748
$ws.Range("D69").Value = 23063
$ws.Range("E69").Value = 15551
$ws.Range("G69").Value = 15
$ws.Range("H69").Value = 818

$ws.Range("A70").Value = "Afganistan"
$ws.Range("B70").Value = 39233
$ws.Range("C70").Value = 6
$ws.Range("D70").Value = 32642
$ws.Range("E70").Value = 5136
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 1455

$ws.Range("A71").Value = "Estado de Palestina"
$ws.Range("B71").Value = 39121
$ws.Range("C71").Value = 418
$ws.Range("D71").Value = 30220
$ws.Range("E71").Value = 8602
$ws.Range("G71").Value = 8
$ws.Range("H71").Value = 299

# Row 127: Surinam - refreshed totals
$ws.Range("B127").Value = 4836
$ws.Range("C127").Value = 1
$ws.Range("D127").Value = 4667
$ws.Range("E127").Value = 67

# Row 147: Guyana - refreshed totals
$ws.Range("B147").Value = 2787
$ws.Range("C147").Value = 15
$ws.Range("D147").Value = 1608
$ws.Range("E147").Value = 1101
$ws.Range("G147").Value = 2
$ws.Range("H147").Value = 78

# Rows 154-155: Burkina Faso's new total overtakes Uruguay, so it moves
# up to row 154 with its refreshed values; Uruguay drops to row 155
# keeping its own (unchanged) figures.
$ws.Range("A154").Value = "Burkina Faso"
$ws.Range("B154").Value = 2028
$ws.Range("C154").Value = 20
$ws.Range("D154").Value = 1279
$ws.Range("E154").Value = 692
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 57

$ws.Range("A155").Value = "Uruguay"
$ws.Range("B155").Value = 2010
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 1755
$ws.Range("E155").Value = 208
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 47

# Row 159: Togo - refreshed totals
$ws.Range("B159").Value = 1749
$ws.Range("C159").Value = 6
$ws.Range("D159").Value = 1336
$ws.Range("E159").Value = 366
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 47

# Row 177: Burundi - refreshed totals
$ws.Range("B177").Value = 502
$ws.Range("C177").Value = 17
$ws.Range("D177").Value = 472
$ws.Range("E177").Value = 29

# Rows 207-208: Santa Lucia and Timor Oriental are tied on every figure,
# but Santa Lucia now sorts first, so the two names simply swap places.
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"
